$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCase")

# --- Fix row 6: the functionality id in column A should be 5, not a duplicate 4 ---
$ws.Cells.Item(6, 1).Value = 5

# --- Turn the two placeholder blank rows (7 and 8) into real test-case rows ---
# Copy row 6's formatting (styles + row height) down into rows 7 and 8 first,
# so the new rows inherit the same borders/fonts/wrap alignment as row 6.
$ws.Range("A6:I6").Copy($ws.Range("A7:I7"))
$ws.Range("A6:I6").Copy($ws.Range("A8:I8"))
$ws.Rows.Item(7).RowHeight = 112.5
$ws.Rows.Item(8).RowHeight = 112.5

# Row 7: check login button with correct email id and blank password
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = "check login button with correct email id and blank password"
$ws.Cells.Item(7, 4).Value = $ws.Cells.Item(6, 4).Value()
$ws.Cells.Item(7, 5).Value = $ws.Cells.Item(6, 5).Value()
$ws.Cells.Item(7, 6).Value = "email id: niravgoti1236@gmail.com`npassword: "
$ws.Cells.Item(7, 7).Value = $ws.Cells.Item(6, 7).Value()
$ws.Cells.Item(7, 8).Value = $ws.Cells.Item(6, 8).Value()
$ws.Cells.Item(7, 9).Value = $ws.Cells.Item(6, 9).Value()

# Row 8: check login button with blank email id and correct password
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 5
$ws.Cells.Item(8, 3).Value = "check login button with blank email id and correct password"
$ws.Cells.Item(8, 4).Value = $ws.Cells.Item(6, 4).Value()
$ws.Cells.Item(8, 5).Value = $ws.Cells.Item(6, 5).Value()
$ws.Cells.Item(8, 6).Value = "email id: `npassword: Jhoncena@15"
$ws.Cells.Item(8, 7).Value = $ws.Cells.Item(6, 7).Value()
$ws.Cells.Item(8, 8).Value = $ws.Cells.Item(6, 8).Value()
$ws.Cells.Item(8, 9).Value = $ws.Cells.Item(6, 9).Value()

# --- Remove the now-superfluous trailing blank row at the bottom of the sheet ---
$ws.Rows.Item(102).Delete()

# --- Update the saved view state to match where the author was working ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A9").Select()
